$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> Vcam1 -> Itgad -> Resolving-Mac) -------------------------
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 17.41485066666667
$ws.Range("H2").Value = 52.24455200000001
$ws.Range("I2").Value = 0.1047285618770465
$ws.Range("J2").Value = 0.1047285618770465
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.366738
$ws.Range("N2").Value = 1.100214
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 6.386687503792001
$ws.Range("R2").Value = 57.48018753412801
$ws.Range("S2").Value = 0.1047285618770465
$ws.Range("T2").Value = 0.1047285618770465

# --- Row 3 (FAPs -> Vcam1 -> Itgad -> Resolving-Mac) ------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 29.68221266666667
$ws.Range("H3").Value = 89.046638
$ws.Range("I3").Value = 0.1785014126970782
$ws.Range("J3").Value = 0.1785014126970782
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.366738
$ws.Range("N3").Value = 1.100214
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 10.885595308948
$ws.Range("R3").Value = 97.970357780532
$ws.Range("S3").Value = 0.1785014126970782
$ws.Range("T3").Value = 0.1785014126970782

# --- Row 4 (MuSCs -> Vcam1 -> Itgad -> Resolving-Mac) -----------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 84.03051233333333
$ws.Range("H4").Value = 252.091537
$ws.Range("I4").Value = 0.5053385113032314
$ws.Range("J4").Value = 0.5053385113032314
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.366738
$ws.Range("N4").Value = 1.100214
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 30.817182032102
$ws.Range("R4").Value = 277.354638288918
$ws.Range("S4").Value = 0.5053385113032314
$ws.Range("T4").Value = 0.5053385113032314

# --- Row 5 (Resolving-Mac -> Vcam1 -> Itgad -> Resolving-Mac) ---------------
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 35.158014
$ws.Range("H5").Value = 105.474042
$ws.Range("I5").Value = 0.2114315141226439
$ws.Range("J5").Value = 0.2114315141226439
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.366738
$ws.Range("N5").Value = 1.100214
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 12.893779738332
$ws.Range("R5").Value = 116.044017644988
$ws.Range("S5").Value = 0.2114315141226439
$ws.Range("T5").Value = 0.2114315141226439

# --- Remove old rows 6-9 (MuSCs/Resolving-Mac combos no longer present) -----
$ws.Rows("6:9").Delete()
